$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column G (hours) from 90 to 120 for all data rows (2 through 160)
$ws.Range("G2:G160").Value = 120
